# Generate Report for Handback
# Updates the handback-status report with refreshed timestamps and a priority
# value change, mirroring a re-run of the report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 2-3
$wsOverview.Range("G2").Value = "2016-11-29 05:31:59"
$wsOverview.Range("G3").Value = "2016-11-29 05:31:59"

# zh-cn sheet: Priority column (E), rows 2-3
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# zh-cn sheet: "Correspond Handoff Datetime" column (H), rows 2-3
$wsZhCn.Range("H2").Value = "2016-11-29 05:31:45"
$wsZhCn.Range("H3").Value = "2016-11-29 05:31:45"

# zh-cn sheet: "Correspond Handback DateTime" column (K), rows 2-3
$wsZhCn.Range("K2").Value = "2016-11-29 05:32:38"
$wsZhCn.Range("K3").Value = "2016-11-29 05:32:38"

# de-de sheet: Priority column (E), rows 2-3
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# de-de sheet: "Correspond Handoff Datetime" column (H), rows 2-3
$wsDeDe.Range("H2").Value = "2016-11-29 05:31:59"
$wsDeDe.Range("H3").Value = "2016-11-29 05:31:59"

# de-de sheet: "Correspond Handback DateTime" column (K), rows 2-3
$wsDeDe.Range("K2").Value = "2016-11-29 05:32:57"
$wsDeDe.Range("K3").Value = "2016-11-29 05:32:57"
